# Add a new worksheet "addNewCustomer" after the existing sheets and
# populate it with a small "add new customer" form, matching the manual
# data-entry session that produced the target workbook.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it becomes sheet3
# (this also makes it the active sheet, which in turn updates
# workbookView/activeTab and moves tabSelected/topLeftCell off the other
# sheets automatically, same as Excel does when you click a new tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "addNewCustomer"

# Give the whole form area a thin black grid border first - matching the
# border style/color already used elsewhere in the workbook - so the cells
# below pick up that existing style (border only) rather than minting a
# brand-new one.
$formRange = $ws.Range("A1:I5")
$formRange.Borders.Color = 0
$formRange.Borders.LineStyle = 1

# Fill in the form, in the same order it was typed originally (label under
# the company name first, then the header row, etc.) so shared-string
# ordering matches.
$ws.Range("A2").Value = "LTD"
$ws.Range("A1").Value = "companyName"
$ws.Range("A3").Value = "Viettel"
$ws.Range("B1").Value = "vatNumber"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "website"
$ws.Range("E1").Value = "groups"
$ws.Range("F1").Value = "address"
$ws.Range("G1").Value = "city"
$ws.Range("H1").Value = "state"
$ws.Range("I1").Value = "zipCode"

$ws.Range("A4").Value = "HBT"
# Leading apostrophe forces these numeric-looking values to be stored as
# text (quote-prefixed), same as the existing "123456" text cell on Sheet1.
$ws.Range("B4").Value = "'667735"
$ws.Range("C4").Value = "'0852741963"
$ws.Range("D4").Value = "hbt.hn.com"
$ws.Range("E4").Value = "khoinghia"
$ws.Range("F4").Value = "Me Linh"
$ws.Range("G4").Value = "Ha Noi"
$ws.Range("H4").Value = "Pass"
$ws.Range("I4").Value = 28386

# Column widths approximating the original auto-fit widths.
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws.Columns.Item(4).ColumnWidth = 11.666666666666666
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 10.0

# Leave the cursor where data entry stopped, one row below the last row.
$ws.Range("A5").Select() | Out-Null
